# Add a "Comment" column (D) to the People data sheet, matching the
# "add comment column to excel data sheet" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("People")

# New header + values for column D.
$ws.Range("D1").Value = "**comment**"
$ws.Range("D2").Value = "aaaa"
$ws.Range("D3").Value = "bbb"
$ws.Range("D4").Value = "ccc"

# Size the new column to fit its contents, like Excel's own best-fit.
$ws.Columns("D").AutoFit() | Out-Null

# Make "People" the active/selected sheet with the new selection used
# while editing the new column.
$ws.Activate() | Out-Null
$ws.Range("E9").Select() | Out-Null
